$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the 'culture_collection' column (Z) entirely - it was re-deleted
# from MIxS per INSDC 2017 confirmation. This shifts everything after
# column Z one column to the left, taking headers and cell comments along.
$ws.Range("Z:Z").Delete()
